$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts them to a numeric value
# (these Price-column entries are intentionally text, e.g. thousands-separated).
$textCells = @("D5", "D6", "D7", "D10", "D11", "D13", "D14", "D17", "D19", "D20", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D33", "D34", "D37", "D38", "D40", "D41", "D42", "D43", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '46.181.27'
$ws.Range('E2').Value = '  +3.80%  '
$ws.Range('D3').Value = '2.451.00'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '320.36'
$ws.Range('E5').Value = '  +2.44%  '
$ws.Range('D6').Value = '104.84'
$ws.Range('E6').Value = '  +3.55%  '
$ws.Range('D7').Value = '0.517'
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +4.83%  '
$ws.Range('D10').Value = '35.93'
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('D11').Value = '0.0806'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('E12').Value = '  -3.22%  '
$ws.Range('D13').Value = '18.26'
$ws.Range('D14').Value = '7.07'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').Value = '2.834.14'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').Value = '2.438.28'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '0.843'
$ws.Range('E17').Value = '  +1.19%  '
$ws.Range('D18').Value = '46.045.89'
$ws.Range('E18').Value = '  +3.72%  '
$ws.Range('D19').Value = '12.58'
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').Value = '6.41'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').Value = '0.0₃0933'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').Value = '71.11'
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '247.24'
$ws.Range('E23').Value = '  +2.43%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '2.36'
$ws.Range('E24').Value = '  +3.75%  '
$ws.Range('D25').Value = '2.50'
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('E26').Value = '  +3.05%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('D29').Value = '9.67'
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('D30').Value = '33.63'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').Value = '49.31'
$ws.Range('E31').Value = '  +1.40%  '
$ws.Range('E32').Value = '  +3.62%  '
$ws.Range('D33').Value = '19.91'
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('D34').Value = '5.34'
$ws.Range('E34').Value = '  +3.35%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').Value = '1.89'
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('D38').Value = '4.51'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('E39').Value = '  +2.36%  '
$ws.Range('D40').Value = '126.89'
$ws.Range('E40').Value = '  +2.66%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '0.111'
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '2.22'
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').Value = '20.83'
$ws.Range('E43').Value = '  -1.63%  '
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('D45').Value = '1.965.15'
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('D46').Value = '2.98'
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('D47').Value = '2.08'
$ws.Range('E47').Value = '  -4.56%  '
$ws.Range('D48').Value = '1.85'
$ws.Range('E48').Value = '  +12.52%  '
$ws.Range('D49').Value = '9.10'
$ws.Range('E49').Value = '  -3.87%  '
$ws.Range('D50').Value = '5.01'
$ws.Range('E50').Value = '  +7.88%  '
$ws.Range('D51').Value = '77.87'
$ws.Range('E51').Value = '  +5.24%  '

Write-Host "Applied 95 cell updates"
